$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws = $wb.Worksheets.Item("Sheet2")

# Rename Sheet2 -> Analysis
$ws.Name = "Analysis"

# Header row: reuse existing shared strings (google_maps_url, trulia_url, zillow_url)
$ws.Range("A1").Value = "google_maps_url"
$ws.Range("B1").Value = "trulia_url"
$ws.Range("C1").Value = "zillow_url"

# Data rows 2..25: HYPERLINK formulas referencing Sheet1 columns E,F,G
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("A$r").Formula = '=HYPERLINK(Sheet1!E' + $r + ', "link")'
    $ws.Range("B$r").Formula = '=HYPERLINK(Sheet1!F' + $r + ', "link")'
    $ws.Range("C$r").Formula = '=HYPERLINK(Sheet1!G' + $r + ', "link")'
}

# Apply the built-in Hyperlink style (adds the Hyperlink font/cellStyle) by
# registering a real hyperlink once, then removing the hyperlink object
# itself (the HYPERLINK() formulas already provide the clickable behaviour) -
# this leaves the style applied to the cells without a <hyperlinks> part.
$ws.Hyperlinks.Add($ws.Range("A2"), "about:blank") | Out-Null
$ws.Range("A2:C25").Style = "Hyperlink"
$ws.Hyperlinks.Delete()

# Column widths to match Sheet1's E/F/G columns
$ws.Columns("A").ColumnWidth = 15.5
$ws.Columns("B").ColumnWidth = 8.333333333333332
$ws.Columns("C").ColumnWidth = 8.833333333333332

# Selection state on the Analysis sheet view, keep Sheet1 as the active tab
$ws.Range("F31").Select() | Out-Null
$ws1.Activate()
